# Apply the edits described by the commit:
# "MAJ + nouveau script 30_application_criteres_lrr + remise à jour des excel
#  concernant les temps de générations et caractère endémique des différents poissons."
#
# Concretely, on sheet "Feuil1" of the active workbook:
#  - Column K's header (K1) is renamed from "duree_generation" to "temps_generation"
#  - Row 14 (Vandoise) A14 flag changes from 0 to 1 (endemic flag)
#  - Many K column values (temps de generation) are updated: some cleared, some
#    updated to new (larger) values reflecting a change of unit / updated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: duree_generation -> temps_generation -------------------
$ws.Range("K1").Value = "temps_generation"

# --- Endemism flag update ---------------------------------------------------
$ws.Range("A14").Value = 1

# --- Column width for column A (species presence column got wider) ---------
$ws.Columns.Item(1).ColumnWidth = 25.5

# --- Clear K values that have no replacement (now unknown/blank) -----------
$clearRows = 4,5,6,8,13,14,15,16,17,20,22,23,25,28,38,42,45
foreach ($r in $clearRows) {
    $ws.Range("K$r").ClearContents()
}

# --- K values updated to new figures ----------------------------------------
$ws.Range("K26").Value = 15
$ws.Range("K27").Value = 12
$ws.Range("K29").Value = 36
$ws.Range("K30").Value = 12
$ws.Range("K31").Value = 15
$ws.Range("K32").Value = 15
$ws.Range("K33").Value = 15
$ws.Range("K34").Value = 10
$ws.Range("K35").Value = 18
$ws.Range("K36").Value = 10
$ws.Range("K37").Value = 10
$ws.Range("K39").Value = 12
$ws.Range("K40").Value = 10
$ws.Range("K41").Value = 15
$ws.Range("K43").Value = 18
$ws.Range("K46").Value = 12
$ws.Range("K47").Value = 15
$ws.Range("K48").Value = 10
$ws.Range("K49").Value = 15
$ws.Range("K50").Value = 10
$ws.Range("K51").Value = 12
$ws.Range("K52").Value = 12

# --- Restore the view / selection state -------------------------------------
$ws.Range("B14").Select()
